$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new columns I and J, matching the formatting of the
# existing header cells (bold font, borders, centered/top alignment).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2 through 40: column I is constant 1, column J mirrors column H
for ($r = 2; $r -le 40; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
